$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.406.90'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.916.94'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4683'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2847'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06803'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '107.46'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.16'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.35%  '
$ws.Range('D12').Value = '1.900.88'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07628'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.176'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6546'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '289.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('D17').Value = '30.418.57'
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007608'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9998'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = '2.158.90'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9990'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.215'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.197'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '21.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.71%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.15'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.266'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.048'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1068'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.373'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.134'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.941'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05029'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7383'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.148'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.735'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02037'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.78%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.686'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.052'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '108.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8728'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.833'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.89%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9995'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '52.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +26.21%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4202'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.167'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.173'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1207'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('B51').Value = 'eCash'
$ws.Range('C51').Value = 'https://coinranking.com/coin/aQx_vW8s1+ecash-xec'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00004439'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +45.66%  '
